$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (but keep formatting) so the shared-string table
# gets rebuilt in the same first-seen order as the source data
# (header, then years 2019..2025).
$ws.Cells.ClearContents()

$headers = @("ano", "ano_obj", "total_customers", "returning_customers", "new_customers", "retention_rate", "new_rate", "returning_rate")

# Final dataset: ano, ano_obj, total_customers, returning_customers, new_customers,
# retention_rate, new_rate, returning_rate
$data = @(
    @("2019", "2019", 2236, 109, 2127, 10.83499005964215, 95.12522361359571, 4.874776386404293),
    @("2020", "2020", 2829, 141, 2688, 6.305903398926654, 95.01590668080594, 4.984093319194061),
    @("2021", "2021", 2493, 184, 2309, 6.504065040650407, 92.61933413557962, 7.380665864420378),
    @("2022", "2022", 2761, 321, 2440, 12.87605294825511, 88.3737776168055, 11.62622238319449),
    @("2023", "2023", 2740, 323, 2417, 11.69865990583122, 88.21167883211679, 11.78832116788321),
    @("2024", "2024", 2445, 304, 2141, 11.09489051094891, 87.56646216768917, 12.43353783231084),
    @("2025", "2025", 688, 130, 558, 5.316973415132924, 81.1046511627907, 18.8953488372093)
)

for ($j = 0; $j -lt $headers.Count; $j++) {
    $col = $j + 1
    $ws.Cells.Item(1, $col).Value = $headers[$j]
}

# Years in columns A/B must stay text (shared-string) cells like the source,
# not get auto-converted to numbers.
$ws.Range("A2:B8").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
